$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. New / modified paragraph & character styles
#    (order matches the target styles.xml: List1, Heading1Char1,
#     List1Char, List3, List4)
# ---------------------------------------------------------------------

# List 1 - paragraph style, based on Heading 1, left-aligned run-in list style
$list1 = $d.Styles.Add("List1", 1)
$list1.NameLocal = "List 1"
$list1.BaseStyle = "Heading1"
$list1.LinkStyle = "List1Char"
$list1.ParagraphFormat.KeepWithNext = $false
$list1.ParagraphFormat.SpaceBefore = 0
$list1.ParagraphFormat.SpaceAfter = 0
$list1.ParagraphFormat.Alignment = 0
$list1.ParagraphFormat.OutlineLevel = 10
$list1.ParagraphFormat.TabStops.Add(184.3, 0)
$list1.Font.Bold = $false
$list1.Font.Size = 12

# Heading 1 Char1 - character style linked with Heading1
$headingChar = $d.Styles.Add("Heading1Char1", 2)
$headingChar.NameLocal = "Heading 1 Char1"
$headingChar.BaseStyle = "DefaultParagraphFont"
$headingChar.LinkStyle = "Heading1"
$headingChar.Font.Bold = $true
$headingChar.Font.Size = 16

$h1 = $d.Styles("Heading1")
$h1.LinkStyle = "Heading1Char1"

# List 1 Char - character style linked with List1
$list1Char = $d.Styles.Add("List1Char", 2)
$list1Char.NameLocal = "List 1 Char"
$list1Char.BaseStyle = "Heading1Char1"
$list1Char.LinkStyle = "List1"
$list1Char.Font.Bold = $false
$list1Char.Font.Size = 12

# re-seat List1's own link now that List1Char exists with the right name
$d.Styles("List1").LinkStyle = "List1Char"

# List 3 - simple indented list paragraph style
$list3 = $d.Styles.Add("List3", 1)
$list3.NameLocal = "List 3"
$list3.BaseStyle = "Normal"
$list3.ParagraphFormat.LeftIndent = 54
$list3.ParagraphFormat.FirstLineIndent = -18
$list3.NoSpaceBetweenParagraphsOfSameStyle = $true

# List 4 - deeper indented list paragraph style, semi-hidden until used
$list4 = $d.Styles.Add("List4", 1)
$list4.NameLocal = "List 4"
$list4.BaseStyle = "Normal"
$list4.UnhideWhenUsed = $true
$list4.ParagraphFormat.LeftIndent = 72
$list4.ParagraphFormat.FirstLineIndent = -18
$list4.NoSpaceBetweenParagraphsOfSameStyle = $true

# List 2 no longer forces Courier New on its run formatting
$list2 = $d.Styles("List2")
$list2.Font.NameAscii = ""
$list2.Font.Name = ""

# ---------------------------------------------------------------------
# 2. Apply the List1 paragraph style to the two lettered sub-paragraphs
# ---------------------------------------------------------------------

$rng1 = $d.Content
$rng1.Find.Execute("   (b) A business ", $true, $false, $false, $false,
                    $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Paragraphs(1).Range.Style = "List1"

$rng2 = $d.Content
$rng2.Find.Execute("   (d)(ii) Requests for approval", $true, $false, $false, $false,
                    $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Paragraphs(1).Range.Style = "List1"
